$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new column before column D. This shifts the existing D:I data
# (Size, Evaluation sub-columns, Time, FPS) one column to the right (E:J)
# and automatically repairs the formulas, the merged-cell ranges and the
# column definitions that reference them.
# ---------------------------------------------------------------------------
$ws.Columns("D").Insert()

# ---------------------------------------------------------------------------
# New column header: "FLOPS (G)" for both tables, mirroring the existing
# "Size (million parameters)" header placement/merge (row3:row4).
# ---------------------------------------------------------------------------
$ws.Range("D3").Value = "FLOPS (G)"
$ws.Range("D4").Value = ""
$ws.Range("D3:D4").Merge()

$ws.Range("D23").Value = "FLOPS (G)"
$ws.Range("D24").Value = ""
$ws.Range("D23:D24").Merge()

# ---------------------------------------------------------------------------
# FLOPS (G) values per model - identical for both benchmark tables, since
# FLOPs depend only on the model architecture, not on the host hardware.
# ---------------------------------------------------------------------------
$flops = @(70, 13.99, 21.76, 13.67, 12.21, 11.07, 7.24, 4.67, 6.39, 6.39, 20.7, 29.58, 15.99, 22.34, 8.76)

for ($i = 0; $i -lt $flops.Length; $i++) {
    $ws.Cells.Item(5 + $i, 4).Value = $flops[$i]
    $ws.Cells.Item(25 + $i, 4).Value = $flops[$i]
}

# ---------------------------------------------------------------------------
# Column D width (narrower than the old "Size" column that is now column E).
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 11.6

# ---------------------------------------------------------------------------
# Restore the view: scroll back to the top-left and move the active
# selection to the new "Size" column on the second table (E32).
# ---------------------------------------------------------------------------
$ws.Range("A1").Select()
$ws.Range("E32").Select()
